# Civil test data: row 3's "occupation/name" column (B3) was stored as a
# blank cell, which downstream code read as the literal string 'NULL'.
# Populate it with an empty-string-safe value ("Ruskin") instead, matching
# the formatting already used by the sibling string cells in the table
# (e.g. D2) rather than leaving it with the sheet's default formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Borrow the formatting of an existing text cell in the same table (D2)
# so the new cell's font matches the rest of the data instead of using
# the worksheet's bare default style.
$ws.Range("D2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats

# Now set the actual value for the previously-blank cell.
$ws.Range("B3").Value = "Ruskin"

# Make sure the font name is explicit (matches the sibling cell's font).
$ws.Range("B3").Font.Name = "Calibri"
